$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1338.7059
$ws.Range("I40").Value = 1349.9286
$ws.Range("K40").Value = 1349.9286
$ws.Range("M40").Value = -1174.9286
$ws.Range("H74").Value = 4710
$ws.Range("J74").Value = 5046.2
$ws.Range("L74").Value = 5046.2
$ws.Range("N74").Value = -6918.2
$ws.Range("H77").Value = 4710
$ws.Range("J77").Value = 5046.2
$ws.Range("L77").Value = 25231
$ws.Range("N77").Value = -34591
$ws.Range("H112").Value = 1706.12
$ws.Range("J112").Value = 1756.7916
$ws.Range("L112").Value = 5270.3748
$ws.Range("N112").Value = -7486.3748
$ws.Range("H137").Value = 1863.7858
$ws.Range("I137").Value = 1449.2941
$ws.Range("K137").Value = 4347.8823
$ws.Range("M137").Value = -1797.8823

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19666.334
$ws.Range("I32").Value = 21094.281
$ws.Range("K32").Value = 21094.281
$ws.Range("M32").Value = -20807.281
$ws.Range("H35").Value = 10700
$ws.Range("I35").Value = 4333.3335
$ws.Range("J35").Value = 29800
$ws.Range("K35").Value = 4333.3335
$ws.Range("L35").Value = 29800
$ws.Range("M35").Value = -3927.3335
$ws.Range("N35").Value = -30612
$ws.Range("H39").Value = 19800
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 19800
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 19800
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -20840
$ws.Range("H61").Value = 6665.579
$ws.Range("I61").Value = 3615.4102
$ws.Range("J61").Value = 13274.277
$ws.Range("K61").Value = 3615.4102
$ws.Range("L61").Value = 13274.277
$ws.Range("M61").Value = -3403.4102
$ws.Range("N61").Value = -13698.277
$ws.Range("H74").Value = 6352.4585
$ws.Range("I74").Value = 2325.85
$ws.Range("J74").Value = 26485.5
$ws.Range("K74").Value = 2325.85
$ws.Range("L74").Value = 26485.5
$ws.Range("M74").Value = -1451.85
$ws.Range("N74").Value = -28233.5
$ws.Range("H77").Value = 6352.4585
$ws.Range("I77").Value = 2325.85
$ws.Range("J77").Value = 26485.5
$ws.Range("K77").Value = 11629.25
$ws.Range("L77").Value = 132427.5
$ws.Range("M77").Value = -7261.25
$ws.Range("N77").Value = -141163.5
$ws.Range("H102").Value = 1613280.8
$ws.Range("I102").Value = 2471590.5
$ws.Range("J102").Value = 3950
$ws.Range("K102").Value = 2471590.5
$ws.Range("L102").Value = 3950
$ws.Range("M102").Value = -2469968.5
$ws.Range("N102").Value = -7194
$ws.Range("H132").Value = 1729.131
$ws.Range("I132").Value = 1318.5469
$ws.Range("K132").Value = 3955.6407
$ws.Range("M132").Value = -1425.6407
$ws.Range("H136").Value = 6665.579
$ws.Range("I136").Value = 3615.4102
$ws.Range("J136").Value = 13274.277
$ws.Range("K136").Value = 10846.2306
$ws.Range("L136").Value = 39822.831
$ws.Range("M136").Value = -8296.230599999999
$ws.Range("N136").Value = -44922.831

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 33334.594
$ws.Range("I134").Value = 1955.16
$ws.Range("J134").Value = 145404
$ws.Range("K134").Value = 5865.48
$ws.Range("L134").Value = 436212
$ws.Range("M134").Value = -3330.48
$ws.Range("N134").Value = -441282

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2212.5
$ws.Range("I31").Value = 1561.6757
$ws.Range("K31").Value = 1561.6757
$ws.Range("M31").Value = -1266.6757
$ws.Range("H34").Value = 2212.5
$ws.Range("I34").Value = 1561.6757
$ws.Range("K34").Value = 1561.6757
$ws.Range("M34").Value = -1359.6757
$ws.Range("H58").Value = 1896007
$ws.Range("I58").Value = 2458289.2
$ws.Range("J58").Value = 4694
$ws.Range("K58").Value = 2458289.2
$ws.Range("L58").Value = 4694
$ws.Range("M58").Value = -2458086.2
$ws.Range("N58").Value = -5100
$ws.Range("H132").Value = 3069.8057
$ws.Range("I132").Value = 3293.3333
$ws.Range("J132").Value = 2526.9524
$ws.Range("K132").Value = 9879.999899999999
$ws.Range("L132").Value = 7580.8572
$ws.Range("M132").Value = -7349.999899999999
$ws.Range("N132").Value = -12640.8572
$ws.Range("H134").Value = 2321.0635
$ws.Range("I134").Value = 1558.2
$ws.Range("J134").Value = 3274.6428
$ws.Range("K134").Value = 4674.6
$ws.Range("L134").Value = 9823.928400000001
$ws.Range("M134").Value = -2139.6
$ws.Range("N134").Value = -14893.9284
$ws.Range("H136").Value = 1896007
$ws.Range("I136").Value = 2458289.2
$ws.Range("J136").Value = 4694
$ws.Range("K136").Value = 7374867.600000001
$ws.Range("L136").Value = 14082
$ws.Range("M136").Value = -7372317.600000001
$ws.Range("N136").Value = -19182

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 470.75
$ws.Range("J7").Value = 383.16666
$ws.Range("L7").Value = 1149.49998
$ws.Range("N7").Value = -1373.49998
$ws.Range("H22").Value = 83334740
$ws.Range("I22").Value = 166666900
$ws.Range("J22").Value = 2575
$ws.Range("K22").Value = 500000700
$ws.Range("L22").Value = 7725
$ws.Range("M22").Value = -500000531
$ws.Range("N22").Value = -8063
$ws.Range("H27").Value = 83334740
$ws.Range("I27").Value = 166666900
$ws.Range("J27").Value = 2575
$ws.Range("K27").Value = 500000700
$ws.Range("L27").Value = 7725
$ws.Range("M27").Value = -500000598
$ws.Range("N27").Value = -7929
$ws.Range("H58").Value = 2960.1177
$ws.Range("J58").Value = 3040.7273
$ws.Range("L58").Value = 9122.1819
$ws.Range("N58").Value = -9378.1819
$ws.Range("H68").Value = 944.1539
$ws.Range("I68").Value = 766.6667
$ws.Range("J68").Value = 997.4
$ws.Range("K68").Value = 2300.0001
$ws.Range("L68").Value = 2992.2
$ws.Range("M68").Value = -1489.0001
$ws.Range("N68").Value = -4614.2
$ws.Range("H71").Value = 944.1539
$ws.Range("I71").Value = 766.6667
$ws.Range("J71").Value = 997.4
$ws.Range("K71").Value = 6900.0003
$ws.Range("L71").Value = 8976.6
$ws.Range("M71").Value = -2844.0003
$ws.Range("N71").Value = -17088.6
$ws.Range("H103").Value = 840.875
$ws.Range("I103").Value = 166.5
$ws.Range("K103").Value = 499.5
$ws.Range("M103").Value = 379.5
$ws.Range("H139").Value = 1469548.4
$ws.Range("I139").Value = 2710893
$ws.Range("J139").Value = 2504.6365
$ws.Range("K139").Value = 8132679
$ws.Range("L139").Value = 7513.9095
$ws.Range("M139").Value = -8127539
$ws.Range("N139").Value = -17793.9095

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 9933.333000000001
$ws.Range("I41").Value = 2500
$ws.Range("J41").Value = 24800
$ws.Range("K41").Value = 2500
$ws.Range("L41").Value = 24800
$ws.Range("M41").Value = -2145
$ws.Range("N41").Value = -25510

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -9828
$ws.Range("N18").ClearContents()
$ws.Range("H61").Value = 731485
$ws.Range("I61").Value = 19898.75
$ws.Range("K61").Value = 19898.75
$ws.Range("M61").Value = -19696.75
$ws.Range("H82").Value = 2642.125
$ws.Range("J82").Value = 2537.4
$ws.Range("L82").Value = 2537.4
$ws.Range("N82").Value = -3259.4
$ws.Range("H85").Value = 2642.125
$ws.Range("J85").Value = 2537.4
$ws.Range("L85").Value = 2537.4
$ws.Range("N85").Value = -5033.4
$ws.Range("H113").Value = 731485
$ws.Range("I113").Value = 19898.75
$ws.Range("K113").Value = 19898.75
$ws.Range("M113").Value = -17728.75
$ws.Range("H132").Value = 5308.3076
$ws.Range("I132").Value = 6085.091
$ws.Range("K132").Value = 18255.273
$ws.Range("M132").Value = -15725.273
$ws.Range("H136").Value = 3328.1943
$ws.Range("I136").Value = 1827.0426
$ws.Range("K136").Value = 5481.1278
$ws.Range("M136").Value = -2931.1278

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5004525
$ws.Range("I5").Value = 3066.8333
$ws.Range("K5").Value = 3066.8333
$ws.Range("M5").Value = -2954.8333
$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1797
$ws.Range("N34").ClearContents()
$ws.Range("H42").Value = 40024.5
$ws.Range("J42").Value = 40024.5
$ws.Range("L42").Value = 40024.5
$ws.Range("N42").Value = -40780.5
$ws.Range("H132").Value = 1790.9736
$ws.Range("I132").Value = 714.75
$ws.Range("J132").Value = 3635.9285
$ws.Range("K132").Value = 2144.25
$ws.Range("L132").Value = 10907.7855
$ws.Range("M132").Value = 385.75
$ws.Range("N132").Value = -15967.7855
$ws.Range("H136").Value = 4786.457
$ws.Range("I136").Value = 3596.9268
$ws.Range("J136").Value = 6468.207
$ws.Range("K136").Value = 10790.7804
$ws.Range("L136").Value = 19404.621
$ws.Range("M136").Value = -8240.7804
